$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33 (VOLTAREN 75MG/3ML 6 AMP.) ---
# Current balance ratio 3:1 -> 3:0
$ws.Range("H33").Value = "3:0"

# Sell price 0.0000 -> 16.3200 (cell is number-formatted text stored as string;
# round-trip the number format so the literal text is kept instead of being
# coerced into a real number)
$cellP33 = $ws.Range("P33")
$cellP33.NumberFormat = "@"
$cellP33.Value = "16.3200"
$cellP33.NumberFormat = "0.00"

# Transactions ratio 0:0 -> 0:1
$ws.Range("Q33").Value = "0:1"

# --- Row 37 (سرنجات 3 سم) ---
$cellP37 = $ws.Range("P37")
$cellP37.NumberFormat = "@"
$cellP37.Value = "4.0000"
$cellP37.NumberFormat = "0.00"

$ws.Range("Q37").Value = "2:0"

# --- Row 38 (سرنجات 5 سم) ---
$cellP38 = $ws.Range("P38")
$cellP38.NumberFormat = "@"
$cellP38.Value = "6.0000"
$cellP38.NumberFormat = "0.00"

$ws.Range("Q38").Value = "2:0"

# --- Row 42 (grand total) ---
$ws.Range("P42").Value = 1479.0550000000001

# --- Row 43 (footer timestamp) ---
$ws.Range("A43").Value = "Wednesday, 17 September, 2025 1:43 PM"
